$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for existing rows 2-31
# from 2024-10-31 (45596) to 2024-11-01 (45597)
for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 3).Value = 45597
}

# Row 31 gains an explicit row height (matching the other data rows)
$ws.Rows.Item(31).RowHeight = 15

# Append a new record in row 32
$row = 32
$ws.Cells.Item($row, 1).Value = "A 47877-2024"
$ws.Cells.Item($row, 2).Value = 45588
$ws.Cells.Item($row, 3).Value = 45597
$ws.Cells.Item($row, 4).Value = "OKÄNT"
$ws.Cells.Item($row, 5).Value = "OKÄNT"
$ws.Cells.Item($row, 7).Value = 1.3
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 18).Value = ""

# Apply the same number formats / styles as the rows above them
$ws.Cells.Item($row, 2).NumberFormat = $ws.Cells.Item(31, 2).NumberFormat
$ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item(31, 3).NumberFormat
$ws.Cells.Item($row, 18).WrapText = $true
